# Update "想去人数" (want-to-go count, column F) figures across all four
# sheets to match freshly scraped bilibili show stats.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 131
$ws.Range("F3").Value = 1311
$ws.Range("F4").Value = 1111
$ws.Range("F5").Value = 996
$ws.Range("F6").Value = 1770
$ws.Range("F8").Value = 1175
$ws.Range("F10").Value = 10
$ws.Range("F12").Value = 279
$ws.Range("F13").Value = 61
$ws.Range("F14").Value = 86
$ws.Range("F15").Value = 664
$ws.Range("F16").Value = 155
$ws.Range("F17").Value = 98
$ws.Range("F21").Value = 132
$ws.Range("F22").Value = 659
$ws.Range("F23").Value = 28
$ws.Range("F28").Value = 309
$ws.Range("F29").Value = 151
$ws.Range("F30").Value = 37
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 250
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 305
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 305
$ws.Range("F3").Value = 131
$ws.Range("F4").Value = 1311
$ws.Range("F5").Value = 1111
$ws.Range("F6").Value = 996
$ws.Range("F7").Value = 1770
$ws.Range("F9").Value = 1175
$ws.Range("F12").Value = 10
$ws.Range("F14").Value = 279
$ws.Range("F15").Value = 61
$ws.Range("F16").Value = 86
$ws.Range("F17").Value = 664
$ws.Range("F18").Value = 155
$ws.Range("F19").Value = 98
$ws.Range("F27").Value = 250
$ws.Range("F28").Value = 250
$ws.Range("F29").Value = 132
$ws.Range("F30").Value = 659
$ws.Range("F31").Value = 28
$ws.Range("F36").Value = 309
$ws.Range("F39").Value = 151
$ws.Range("F40").Value = 37
